# ---------------------------------------------------------------------------
# Update with Correct Forecast output
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ---- 1. Rename the original sheet and wire up the new analysis sheets ----
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sales vs PO"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

# ---------------------------------------------------------------------------
# 2. "Sales vs PO" sheet: insert a new "Order Week" column (C) that carries
#    the original weekly date, shift the "ds" dates forward by 6 days, and
#    zero out the PO_Requested_Qty column (now D) - those figures moved to
#    the new "Weekly Growth" sheet.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).Insert()

# copy date formatting from column A onto the newly inserted column C
$ws1.Range("A1:A18").Copy()
$ws1.Range("C1:C18").PasteSpecial(-4122)

$ws1.Range("C1").Value = "Order Week"
$ws1.Range("D1").Value = "PO_Requested_Qty"

$dsVals = @(45543,45550,45557,45564,45571,45578,45585,45592,45599,45606,45613,45620,45627,45634,45641,45648,45655)
$yVals = @(0,0,0,14,56,38,39,42,100,36,8,46,19,135,119,43,52)
$orderWeekVals = @(45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)

for ($i = 0; $i -lt $dsVals.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $dsVals[$i]
    $ws1.Cells.Item($row, 2).Value = $yVals[$i]
    $ws1.Cells.Item($row, 3).Value = $orderWeekVals[$i]
    $ws1.Cells.Item($row, 4).Value = 0
}

# ---------------------------------------------------------------------------
# 3. "Weekly Growth" sheet: the PO quantities per order week together with
#    the week-over-week growth percentage.
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$growthDs = @(45544,45572,45586,45593,45614)
$growthQty = @(600,230,450,480,1560)
$growthPct = @(0,-61.66666666666667,95.65217391304348,6.666666666666665,225)

for ($i = 0; $i -lt $growthDs.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $growthDs[$i]
    $ws2.Cells.Item($row, 2).Value = $growthQty[$i]
    $ws2.Cells.Item($row, 3).Value = $growthPct[$i]
}

# copy the header style (bold + border) and date number formatting (style)
# used on "Sales vs PO" onto the "Weekly Growth" sheet
$ws1.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. "Volume Insights" sheet: summary stats over the PO quantities.
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 3320
$ws3.Range("B2").Value = 664
$ws3.Range("C2").Value = 1560
$ws3.Range("D2").Value = 230

$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. "Prediction Info" sheet: the next-week PO quantity forecast.
# ---------------------------------------------------------------------------
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 1315

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. Keep "Sales vs PO" as the active/selected sheet (as in the original
#    workbook) since adding the extra sheets shifts focus to the last one.
# ---------------------------------------------------------------------------
$ws1.Activate()
